$wb = $excel.ActiveWorkbook

# Rename the sheet "excel" -> "excel (1)"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "excel (1)"

# Rename the defined name "excel" -> "excel__1" (its RefersTo is
# automatically re-pointed at the renamed sheet by the rename above)
$n = $wb.Names.Item(1)
$n.Name = "excel__1"

# Update the "report generated" timestamp text (stored as a shared string)
$cell = $ws.Cells.Item(24, 1)
$cell.Value = "Report generated at 9:51 PM on Dec 6, 2018"
